# Update expected value in Sheet1!B2 from "y" to "n" (float support now
# expected to be handled, so the prior "y" assertion toggle is no longer
# applicable and the test now expects "n").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make sure Sheet1 is the active sheet/view (it was already tabSelected).
$ws.Activate()

$ws.Range("B2").Value = "n"

# Reset the view: scroll back to the top-left (A1) and move the active
# selection to B6, matching the post-edit cursor position.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B6").Select()
